# Swap the content of rows 3 and 4 for the columns that differ between
# them (A, B, D, E, F, G, H, Q, R, AC), leaving the columns that are
# identical between the two rows untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr3 = "$col`3"
    $addr4 = "$col`4"
    $v3 = $ws.Range($addr3).Value2
    $v4 = $ws.Range($addr4).Value2
    $ws.Range($addr3).Value2 = $v4
    $ws.Range($addr4).Value2 = $v3
}

# The "Publik kommentar" text moves from row 3 to row 4.
$comment = $ws.Range("AC3").Value2
$ws.Range("AC3").Value2 = $null
$ws.Range("AC4").Value2 = $comment
